# 17.1.1 — add the 2021 / 2022 columns (R, S) to the government-revenue
# table and refresh the already-published 2019 / 2020 (P, Q) figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Copy-CellFormat($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats) | Out-Null
}

# ---- Row 4 : years header -------------------------------------------------
Copy-CellFormat "Q4" "R4"
Copy-CellFormat "Q4" "S4"
$ws.Range("R4").Value = 2021
$ws.Range("S4").Value = 2022

# ---- Row 5 : Revenues, total -----------------------------------------------
Copy-CellFormat "Q5" "R5"
Copy-CellFormat "Q5" "S5"
$ws.Range("P5").Value = 25.6
$ws.Range("Q5").Value = 23.8
$ws.Range("R5").Value = 26.8
$ws.Range("S5").Value = 26.8

# ---- Row 6 : Tax revenues ---------------------------------------------------
Copy-CellFormat "Q6" "R6"
Copy-CellFormat "Q6" "S6"
$ws.Range("P6").Value = 18.6
$ws.Range("Q6").Value = 16.7
$ws.Range("R6").Value = 19.3
$ws.Range("S6").Value = 19.3

# ---- Row 7 : Contributions / deductions for social needs -------------------
Copy-CellFormat "Q7" "R7"
Copy-CellFormat "Q7" "S7"
$ws.Range("R7").Value = "-"
$ws.Range("S7").Value = "-"

# ---- Row 8 : Received official transfers ------------------------------------
Copy-CellFormat "Q8" "R8"
Copy-CellFormat "Q8" "S8"
$ws.Range("P8").Value = 2.1
$ws.Range("Q8").Value = 1.8
$ws.Range("R8").Value = 1.8
$ws.Range("S8").Value = 1.8

# ---- Row 9 : Non-tax revenues ------------------------------------------------
Copy-CellFormat "Q9" "R9"
Copy-CellFormat "Q9" "S9"
$ws.Range("P9").Value = 4.9
$ws.Range("Q9").Value = 5.2
$ws.Range("R9").Value = 5.7
$ws.Range("S9").Value = 5.7

# ---- Row 10 : Revenues from the sale of non-financial assets ----------------
Copy-CellFormat "Q10" "R10"
Copy-CellFormat "Q10" "S10"
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 0

$excel.CutCopyMode = 0

# ---- view bookkeeping --------------------------------------------------------
$ws.Range("T3").Select() | Out-Null
